$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("D222").Value = 123
$ws.Range("F222").Value = 121
$ws.Range("M222").Formula = "=ABS(D222-F222)"
$ws.Range("M222").NumberFormat = $ws.Range("M221").NumberFormat
Write-Host "M221 numfmt: $($ws.Range("M221").NumberFormat)"
